$d = $word.ActiveDocument

# Locate the paragraph that ends with "Once you've changed follow the next step in your assignment task."
$anchor = $d.Paragraphs.Item(5)
$anchorRange = $anchor.Range

# Insert four new paragraphs after it (each inherits the bordered paragraph
# formatting from the anchor paragraph automatically).
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()

# The first three new paragraphs should stay empty (no run at all). Typing a
# placeholder character and then removing it via Find/Replace drops the
# empty run entirely, matching a paragraph that never received a run.
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertBefore("X")
$p6.Range.Find.Execute("X", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertBefore("X")
$p7.Range.Find.Execute("X", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

$p8 = $d.Paragraphs.Item(8)
$p8.Range.InsertBefore("X")
$p8.Range.Find.Execute("X", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# The fourth new paragraph gets real text plus a _GoBack bookmark marking the
# last edit position (what Word stamps automatically on save).
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Text = "Add new content to the file for testing"
$d.Bookmarks.Add("_GoBack", $p9.Range)
